# Narrative.docx restructuring script
# Implements: new title paragraph, rewritten Synopsis, new Backstory / Audio
# logs sections, relocated "Player character" section, updated Ship details,
# and removal of the AI/"Characters in audio log" section.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Remove the trailing "AI: Arin" + empty paragraph + "Characters in
# audio log" block (old paragraphs 30-35).
# ---------------------------------------------------------------------------
$pAiStart = $d.Paragraphs.Item(30)
$pAiEnd = $d.Paragraphs.Item(35)
$removeRange = $d.Range($pAiStart.Range.Start, $pAiEnd.Range.End)
$removeRange.Delete()

# ---------------------------------------------------------------------------
# Step 2: Update the "Ship" section - add "Name: " / "Type: " labels and
# extend the vessel description.
# ---------------------------------------------------------------------------
$pShipName = $d.Paragraphs.Item(26)
$nameStart = $pShipName.Range.Start
$d.Range($nameStart, $nameStart).InsertBefore("Name: ")

$pShipType = $d.Paragraphs.Item(27)
$pShipType.Range.Find.Execute(
    "Orion-class science vessel", $true, $false, $false, $false, $false,
    $true, 1, $false, "Type: Orion-class science and exploration vessel", 2)

# ---------------------------------------------------------------------------
# Step 3: Replace the old Synopsis body / spoiler / goal / "###" block
# (old paragraphs 2-14) with the new Synopsis text, Backstory section and
# Audio logs section.
# ---------------------------------------------------------------------------
$pBlockStart = $d.Paragraphs.Item(2)
$pBlockEnd = $d.Paragraphs.Item(14)
$blockRange = $d.Range($pBlockStart.Range.Start, $pBlockEnd.Range.End)
$blockRange.Text = "placeholder"

# Paragraph 2: new Synopsis paragraph text.
$d.Paragraphs.Item(2).Range.Text = "You’re stranded on a strange and beautiful chain of islands in space. To send a distress call " + [char]0x2014 + " which will help you get home " + [char]0x2014 + " you will explore the islands, gather parts, and repair a communications tower. But before you leave, a poignant transmission will help you deal with a painful part of your past and transition from reconstruction to acceptance."

# Paragraph 3 stays empty (already blank after the big Text= assignment).

# Paragraph 4: "Backstory" heading.
$pBackstoryHeading = $d.Paragraphs.Item(4)
$pBackstoryHeading.Range.Text = "Backstory"
$pBackstoryHeading.Range.Font.Bold = $true

# Paragraph 5: Backstory body text.
$d.Paragraphs.Item(5).Range.Text = "Terry’s four-year-old, Alex, died from a disease 10 years ago. The island represents Terry’s emotional journey from reconstruction to grieving. The art style represents a combination of the literal and surreal nature of Terry’s journey."

# Paragraph 6 stays empty.

# Paragraph 7: "Audio logs" heading.
$pAudioHeading = $d.Paragraphs.Item(7)
$pAudioHeading.Range.Text = "Audio logs"
$pAudioHeading.Range.Font.Bold = $true

# Paragraph 8: Audio logs body text.
$d.Paragraphs.Item(8).Range.Text = "The narrative unfolds through a series of audio logs, which reveal the backstory and the ultimate transmission."

# Paragraph 9 stays empty.

# Paragraphs 10-14 are left-over placeholders from the original 13-paragraph
# block (we only needed 8: 2-9); remove them.
$pExtraStart = $d.Paragraphs.Item(10)
$pExtraEnd = $d.Paragraphs.Item(14)
$extraRange = $d.Range($pExtraStart.Range.Start, $pExtraEnd.Range.End)
$extraRange.Delete()

# ---------------------------------------------------------------------------
# Step 4: Insert the new bold, centered title paragraph at the very top of
# the document.
# ---------------------------------------------------------------------------
$pFirst = $d.Paragraphs.Item(1)
$pFirst.Range.InsertParagraphBefore()
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Alignment = 1
$titlePara.Range.Text = [char]0x201C + "Acceptance" + [char]0x201D + " narrative"
$titlePara.Range.Font.Bold = $true
$titlePara.Range.Font.Size = 14
$titlePara.Range.Font.SizeBi = 14
